$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 6
    3  = 3
    4  = 5
    5  = 6
    6  = 8
    7  = 3
    8  = 6
    9  = 6
    10 = 3
    11 = 8
    12 = 7
    13 = 7
    14 = 3
    15 = 6
    16 = 7
    17 = 8
    18 = 5
    19 = 8
    20 = 5
    21 = 8
    22 = 3
    23 = 5
    24 = 8
    25 = 11
    26 = 6
    27 = 8
    28 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
